{"js": "// Update the \"Overall\" row's endline prevalence figures (SDB prevalence\n// table) to reflect the new 58-entry analysis:\n//   SDBendline:    17.72(9.3-26.14)    -> 24.14(13.13-35.15)\n//   no SDBendline: 26.58(16.84-36.32)  -> 36.21(23.84-48.58)\n//\n// The same \"26.58(16.84-36.32)\" string also appears in the \"no\n// SDBBaseline\" column of the very same row, so we locate the specific\n// table cells (by header text / row label) rather than doing a blind\n// document-wide text replace, to avoid touching the baseline value.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"values\");\nawait context.sync();\n\nconst values = table.values;\nconst header = values[0];\n\n// Locate the \"SDBendline\" and \"no SDBendline\" columns from the header row.\nlet sdbEndlineCol = -1;\nlet noSdbEndlineCol = -1;\nfor (let c = 0; c < header.length; c++) {\n  const text = header[c];\n  if (text.indexOf(\"no SDB\") === 0 && text.indexOf(\"endline\") > 0) {\n    noSdbEndlineCol = c;\n  } else if (text.indexOf(\"SDB\") === 0 && text.indexOf(\"endline\") > 0) {\n    sdbEndlineCol = c;\n  }\n}\n\n// Locate the \"Overall\" row from the first column.\nlet overallRow = -1;\nfor (let r = 0; r < values.length; r++) {\n  if (values[r][0] === \"Overall\") {\n    overallRow = r;\n    break;\n  }\n}\n\n// Replace each value, scoped tightly to its own cell so the identical\n// string living elsewhere in the same row/table is left untouched.\nconst endlineCell = table.getCell(overallRow, sdbEndlineCol);\nconst noEndlineCell = table.getCell(overallRow, noSdbEndlineCol);\n\nconst endlineMatches = endlineCell.body.search(\"17.72(9.3-26.14)\", { matchCase: true });\nendlineMatches.load(\"items\");\nconst noEndlineMatches = noEndlineCell.body.search(\"26.58(16.84-36.32)\", { matchCase: true });\nnoEndlineMatches.load(\"items\");\nawait context.sync();\n\nendlineMatches.items[0].insertText(\"24.14(13.13-35.15)\", Word.InsertLocation.replace);\nnoEndlineMatches.items[0].insertText(\"36.21(23.84-48.58)\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Update the \"Overall\" row's endline prevalence figures (SDB prevalence\n# table) to reflect the new 58-entry analysis:\n#   SDBendline:    17.72(9.3-26.14)    -> 24.14(13.13-35.15)\n#   no SDBendline: 26.58(16.84-36.32)  -> 36.21(23.84-48.58)\n#\n# The same \"26.58(16.84-36.32)\" string also appears in the \"no\n# SDBBaseline\" column of the very same row, so we locate the specific\n# table cells (by header text / row label) rather than doing a blind\n# document-wide replace, to avoid touching the baseline value.\n\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\n# Locate the \"SDBendline\" and \"no SDBendline\" columns from the header row.\n$sdbEndlineCol = $null\n$noSdbEndlineCol = $null\nfor ($c = 1; $c -le $tbl.Columns.Count; $c++) {\n    $headerText = $tbl.Cell(1, $c).Range.Text\n    if ($headerText -like \"no SDB*endline*\") {\n        $noSdbEndlineCol = $c\n    } elseif ($headerText -like \"SDB*endline*\") {\n        $sdbEndlineCol = $c\n    }\n}\n\n# Locate the \"Overall\" row from the first column.\n$overallRow = $null\nfor ($r = 1; $r -le $tbl.Rows.Count; $r++) {\n    $label = $tbl.Cell($r, 1).Range.Text\n    $label = $label.Substring(0, $label.Length - 2)\n    if ($label -eq \"Overall\") {\n        $overallRow = $r\n        break\n    }\n}\n\n# Replace the SDBendline value, scoped tightly to that one cell so the\n# identical string elsewhere in the document/table is left untouched.\n$cellA = $tbl.Cell($overallRow, $sdbEndlineCol)\n$rngA = $d.Range($cellA.Range.Start, $cellA.Range.End)\n$findA = $rngA.Find\n$findA.Text = \"17.72(9.3-26.14)\"\n$findA.Replacement.Text = \"24.14(13.13-35.15)\"\n$findA.Execute([ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, 1) | Out-Null\n\n# Replace the no-SDBendline value, scoped to its own cell (column 7),\n# not the duplicate \"26.58(16.84-36.32)\" sitting in the baseline column.\n$cellB = $tbl.Cell($overallRow, $noSdbEndlineCol)\n$rngB = $d.Range($cellB.Range.Start, $cellB.Range.End)\n$findB = $rngB.Find\n$findB.Text = \"26.58(16.84-36.32)\"\n$findB.Replacement.Text = \"36.21(23.84-48.58)\"\n$findB.Execute([ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, 1) | Out-Null\n"}
